$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number + week-covering dates) ---
$ws.Range("A8").Value = "Volume 32   Number  26"
$ws.Range("C9").Value = "Report Covering the Week  6/23/2025  Through  6/29/2025"

# --- Weekly crime statistics table updates ---
# Row 15
$ws.Range("D15").Value = 2
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = -75
$ws.Range("J15").Value = 7
$ws.Range("K15").Value = -14.285714285714

# Row 16
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = -6.25
$ws.Range("I16").Value = 60
$ws.Range("J16").Value = 69
$ws.Range("K16").Value = -13.043478260869
$ws.Range("L16").Value = -13.043478260869
$ws.Range("M16").Value = -23.076923076923
$ws.Range("N16").Value = -85

# Row 17
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = -44.444444444444
$ws.Range("F17").Value = 19
$ws.Range("G17").Value = 31
$ws.Range("H17").Value = -38.709677419354
$ws.Range("I17").Value = 115
$ws.Range("J17").Value = 137
$ws.Range("K17").Value = -16.058394160583
$ws.Range("L17").Value = 7.476635514018
$ws.Range("M17").Value = 112.962962962963
$ws.Range("N17").Value = -24.836601307189

# Row 18
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("I18").Value = 74
$ws.Range("J18").Value = 78
$ws.Range("K18").Value = -5.128205128205
$ws.Range("L18").Value = -33.928571428571
$ws.Range("M18").Value = -47.142857142857
$ws.Range("N18").Value = -89.606741573033

# Row 19
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = -42.857142857142
$ws.Range("G19").Value = 52
$ws.Range("H19").Value = -5.769230769230
$ws.Range("I19").Value = 268
$ws.Range("J19").Value = 312
$ws.Range("K19").Value = -14.102564102564
$ws.Range("L19").Value = -16.510903426791
$ws.Range("M19").Value = 29.468599033816
$ws.Range("N19").Value = -26.975476839237

# Row 20
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 15
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 88
$ws.Range("J20").Value = 88
$ws.Range("L20").Value = 29.411764705882
$ws.Range("M20").Value = 3.529411764705
$ws.Range("N20").Value = -89.908256880733

# Row 21
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 32
$ws.Range("E21").Value = -37.5
$ws.Range("F21").Value = 111
$ws.Range("G21").Value = 132
$ws.Range("H21").Value = -15.909090909090
$ws.Range("I21").Value = 611
$ws.Range("J21").Value = 693
$ws.Range("K21").Value = -11.832611832611
$ws.Range("L21").Value = -10.932944606414
$ws.Range("M21").Value = 6.445993031358
$ws.Range("N21").Value = -75.734710087370

# Row 22
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = -100
$ws.Range("G22").Value = 7
$ws.Range("H22").Value = -42.857142857142
$ws.Range("J22").Value = 15
$ws.Range("K22").Value = -20
$ws.Range("L22").Value = 100
$ws.Range("M22").Value = -29.411764705882

# Row 24
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 29
$ws.Range("E24").Value = -20.689655172413
$ws.Range("F24").Value = 99
$ws.Range("G24").Value = 115
$ws.Range("H24").Value = -13.913043478260
$ws.Range("I24").Value = 811
$ws.Range("J24").Value = 857
$ws.Range("K24").Value = -5.367561260210
$ws.Range("L24").Value = -10.485651214128
$ws.Range("M24").Value = 75.161987041036

# Row 25
$ws.Range("D25").Value = 16
$ws.Range("E25").Value = -31.25
$ws.Range("G25").Value = 66
$ws.Range("H25").Value = -28.787878787878
$ws.Range("I25").Value = 450
$ws.Range("J25").Value = 488
$ws.Range("K25").Value = -7.786885245901
$ws.Range("L25").Value = -13.957934990439

# Row 26
$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 15
$ws.Range("E26").Value = -40
$ws.Range("F26").Value = 41
$ws.Range("G26").Value = 58
$ws.Range("H26").Value = -29.310344827586
$ws.Range("I26").Value = 319
$ws.Range("J26").Value = 320
$ws.Range("K26").Value = -0.3125
$ws.Range("L26").Value = 6.688963210702
$ws.Range("M26").Value = 34.033613445378

# Row 27
$ws.Range("D27").Value = 2
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -75
$ws.Range("J27").Value = 11
$ws.Range("K27").Value = -36.363636363636

# Row 28
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 9
$ws.Range("H28").Value = -55.555555555555
$ws.Range("J28").Value = 24
$ws.Range("K28").Value = -8.333333333333
$ws.Range("L28").Value = -15.384615384615

# Row 31
$ws.Range("G31").Value = 2
$ws.Range("H31").Value = -50

# --- Cells that become text "0" (matching style of neighboring text/N-A cells) ---
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("C23").Copy()
$ws.Range("C22").PasteSpecial(-4122)

$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "0"
$ws.Range("C27").Copy()
$ws.Range("C28").PasteSpecial(-4122)

$excel.CutCopyMode = 0
